# Almost integrated section 5's workflow:
# - Roll the "Selection 1 / Selection 2" review dates forward a month's worth of days
#   (Selection 1: 01-09-2020 -> 01-22-2020, Selection 2: 01-22-2020 -> 01-24-2020)
# - Shift the change-list evidence names that go with those selections, dropping the
#   oldest set and adding a new (still partially filled in) set for Selection 2
# - Remove the stray leftover cell that had drifted into column H on row 35 and let the
#   rest of the "step 4" paragraph collapse back up into place
# - Mirror the same evidence-name shuffle on the Data tab

$wb = $excel.ActiveWorkbook
$report = $wb.Sheets("Report")
$data = $wb.Sheets("Data")

# --- Step 2 block: Selection 1/2 (date) + Emails Received ---
$report.Range("E16").Value = 43852   # Selection 1 (date): 1/9/2020 -> 1/22/2020
$report.Range("E17").Value = 43854   # Selection 2 (date): 1/22/2020 -> 1/24/2020
$report.Range("H17").Value = 3       # Emails Received for Selection 2: 4 -> 3

# --- Step 3 block: Selection 1 (Date) / Total Changes list ---
$report.Range("D28").Value = 43852   # Selection 1 (Date): 1/9/2020 -> 1/22/2020
$report.Range("H28").Value = "magic435"
$report.Range("H29").Value = "magic_qq_23455"
$report.Range("H30").Value = "test_2348q"
$report.Range("H31").Value = "test234234"

# --- Step 3 block: Selection 2 (Date) / Total Changes list ---
$report.Range("D32").Value = 43854   # Selection 2 (Date): 1/22/2020 -> 1/24/2020
$report.Range("G32").Value = 3       # Total Changes for Selection 2: 4 -> 3
$report.Range("H32").Value = "magic_iq23"
$report.Range("H33").Value = "testqcl12"
$report.Range("H34").Value = "oiuer3298"

# Row 35 only held a stray leftover value in column H (left over from the old
# Selection 2 list). Delete the whole row so the "4. From the Findings Tracker..."
# paragraph (previously starting on row 36) shifts back up to row 35, and
# everything beneath it shifts up by one row too.
$report.Rows("35:35").Delete()

# Reset the view back to the top of the sheet and leave the cursor on G31.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$report.Range("G31").Select()

# --- Data tab: shift the evidence file names the same way ---
$data.Range("A1").Value = "C:\Users\karth\OneDrive\Documents\UiPath\Robotic_Process_Automation\FOLDER\Mainfolder\Daily_Change_Monitoring\1Jan2020\01-22-2020\*"
$data.Range("B1").Value = "C:\Users\karth\OneDrive\Documents\UiPath\Robotic_Process_Automation\FOLDER\Mainfolder\Daily_Change_Monitoring\1Jan2020\01-24-2020\*"

$data.Range("A2").Value = "CHANGES - SOX Audit Report for magic435.txt_07.01.73"
$data.Range("B2").Value = "CHANGES - SOX Audit Report for magic_iq23.txt_07.01.73"

$data.Range("A3").Value = "CHANGES - SOX Audit Report for magic_qq_23455.txt_07.01.73"
$data.Range("B3").Value = "CHANGES - SOX Audit Report for oiuer3298.txt_07.01.73"

$data.Range("A4").Value = "CHANGES - SOX Audit Report for test234234.txt_07.01.73"
$data.Range("B4").Value = "CHANGES - SOX Audit Report for testqcl12.txt_07.01.73"

$data.Range("A5").Value = "CHANGES - SOX Audit Report for test_2348q.txt_07.01.73"
